$d = $word.ActiveDocument
$wNS = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Find-ParaByPrefix($prefix) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

function Replace-ParaXml($prefix, $newXml) {
    $p = Find-ParaByPrefix($prefix)
    if ($null -eq $p) {
        Write-Output ("NOT FOUND: " + $prefix)
        return
    }
    $r = $p.Range
    $r.InsertXML($newXml)
}

# --- 1) Title paragraph: drop trailing <w:br/> run, then insert the new
#        "Packages used in this exercise:" paragraph (with the library()
#        lines) right after it, ending with the break that used to close
#        the title paragraph. ---
$titleXml = "<w:p $wNS><w:r><w:t xml:space='preserve'>Exercise </w:t></w:r><w:r><w:t>3</w:t></w:r><w:r><w:t xml:space='preserve'>: </w:t></w:r><w:r><w:t>Size structure and weight-length relationship</w:t></w:r></w:p>" + `
            "<w:p $wNS><w:r><w:t>Packages used in this exercise:</w:t></w:r><w:r><w:br/></w:r><w:proofErr w:type='gramStart'/><w:r><w:t>library(</w:t></w:r><w:proofErr w:type='gramEnd'/><w:r><w:t>FSA)</w:t></w:r><w:r><w:br/><w:t>library(ggplot2)</w:t></w:r><w:r><w:br/><w:t>library(</w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>tidyr</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>)</w:t></w:r><w:r><w:br/><w:t>library(</w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>dplyr</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>)</w:t></w:r><w:r><w:br/></w:r></w:p>"
Replace-ParaXml "Exercise " $titleXml

# --- 2) "Generate a vector..." bullet: split the BluegillLM / FSAdata
#        mentions with spell-check proofErr wraps, then add a trailing
#        empty ListParagraph. ---
$p1Xml = "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Generate a vector of percentage of fish in 10mm length intervals using </w:t></w:r><w:r><w:t>the “</w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>BluegillLM</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space='preserve'> data set from the </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>FSAdata</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> package.</w:t></w:r></w:p>" + `
         "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"
Replace-ParaXml "Generate a vector" $p1Xml

# --- 3) "Create a length frequency histogram..." bullet ---
$p2Xml = "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Create a length frequency histogram of total length (</w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>tl</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> column) using the “</w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>BluegillLM</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space='preserve'> data set from the </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>FSAdata</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> package.</w:t></w:r><w:r><w:t xml:space='preserve'> Hint: use the vector created in question 1 to determine range of breaks for the histogram</w:t></w:r></w:p>" + `
         "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"
Replace-ParaXml "Create a length frequency histogram" $p2Xml

# --- 4) "Create a frequency table..." bullet ---
$p3Xml = "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Create a frequency table of Bluegill PSD size groups using </w:t></w:r><w:r><w:t>the “</w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>BluegillLM</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space='preserve'> data set from the </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>FSAdata</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> package. </w:t></w:r></w:p>" + `
         "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"
Replace-ParaXml "Create a frequency table" $p3Xml

# --- 5) "Calculate the PSD-Q and PSD-P..." bullet ---
$p4Xml = "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Calculate the PSD-Q and PSD-P of Bluegill using the </w:t></w:r><w:r><w:t>“</w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>BluegillLM</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space='preserve'> dataset from the </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>FSAdata</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> package.</w:t></w:r></w:p>" + `
         "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"
Replace-ParaXml "Calculate the PSD-Q" $p4Xml

# --- 6) "Create a scatterplot..." bullet ---
$p5Xml = "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Create a scatterplot of Bluegill total length (x-axis) and weight (y-axis) on the natural log scale using </w:t></w:r><w:r><w:t>the “</w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>BluegillLM</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space='preserve'> data set from the </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>FSAdata</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> package.</w:t></w:r></w:p>" + `
         "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"
Replace-ParaXml "Create a scatterplot" $p5Xml

# --- 7) "Estimate coefficients..." bullet (last bullet, no trailing
#        empty paragraph is added after this one). ---
$p6Xml = "<w:p $wNS><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr><w:r><w:t xml:space='preserve'>Estimate coefficients of a weight-length model using </w:t></w:r><w:r><w:t>the “</w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>BluegillLM</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t>”</w:t></w:r><w:r><w:t xml:space='preserve'> data set from the </w:t></w:r><w:proofErr w:type='spellStart'/><w:r><w:t>FSAdata</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> package. Also determine the 95% confidence intervals for the intercept and slope.</w:t></w:r></w:p>"
Replace-ParaXml "Estimate coefficients" $p6Xml

Write-Output "Done"
